$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Configuration_template")
$ws2 = $wb.Worksheets.Item("Lable")

# Select B9 on Configuration_template sheet
$ws1.Activate()
$ws1.Range("B9").Select()

# Sort the Lable sheet data A4:A12 ascending
$rng = $ws2.Range("A4:A12")
$rng.Sort($ws2.Range("A4"), 1)
